$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the ShipmentTracking values (column P) for rows 3, 4 and 5 with
# the new FedEx tracking numbers recorded on 15th June 2022.
# These tracking numbers are stored as text, so the cell format is set
# to Text first - otherwise the all-digit strings would be auto-converted
# to numbers.
$ws.Range("P3:P5").NumberFormat = "@"
$ws.Range("P3").Value = "320018812762"
$ws.Range("P4").Value = "320018812800"
$ws.Range("P5").Value = "320018812810"
